$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = [double]"1"
$ws.Range("F2").Value2 = [double]"0.3333333333333333"
$ws.Range("G2").Value2 = [double]"0.03034166666666667"
$ws.Range("H2").Value2 = [double]"0.09102499999999999"
$ws.Range("I2").Value2 = [double]"0.001862306623420098"
$ws.Range("J2").Value2 = [double]"0.001862306623420098"
$ws.Range("M2").Value2 = [double]"0.07271233333333334"
$ws.Range("N2").Value2 = [double]"0.218137"
$ws.Range("O2").Value2 = [double]"0.004171225362010892"
$ws.Range("P2").Value2 = [double]"0.004171225362010893"
$ws.Range("Q2").Value2 = [double]"0.002206213380555556"
$ws.Range("R2").Value2 = [double]"0.019855920425"
$ws.Range("S2").Value2 = [double]"7.768100619450781E-06"
$ws.Range("T2").Value2 = [double]"7.768100619450781E-06"
$ws.Range("E3").Value2 = [double]"1"
$ws.Range("F3").Value2 = [double]"0.3333333333333333"
$ws.Range("G3").Value2 = [double]"0.03034166666666667"
$ws.Range("H3").Value2 = [double]"0.09102499999999999"
$ws.Range("I3").Value2 = [double]"0.001862306623420098"
$ws.Range("J3").Value2 = [double]"0.001862306623420098"
$ws.Range("O3").Value2 = [double]"0.5387060579248023"
$ws.Range("P3").Value2 = [double]"0.5387060579248023"
$ws.Range("Q3").Value2 = [double]"0.2849283867527778"
$ws.Range("R3").Value2 = [double]"2.564355480775"
$ws.Range("S3").Value2 = [double]"0.00100323585974989"
$ws.Range("T3").Value2 = [double]"0.00100323585974989"
$ws.Range("E4").Value2 = [double]"1"
$ws.Range("F4").Value2 = [double]"0.3333333333333333"
$ws.Range("G4").Value2 = [double]"0.03034166666666667"
$ws.Range("H4").Value2 = [double]"0.09102499999999999"
$ws.Range("I4").Value2 = [double]"0.001862306623420098"
$ws.Range("J4").Value2 = [double]"0.001862306623420098"
$ws.Range("O4").Value2 = [double]"0.4571227167131868"
$ws.Range("P4").Value2 = [double]"0.4571227167131868"
$ws.Range("Q4").Value2 = [double]"0.2417779349333334"
$ws.Range("R4").Value2 = [double]"2.1760014144"
$ws.Range("S4").Value2 = [double]"0.0008513026630507568"
$ws.Range("T4").Value2 = [double]"0.0008513026630507567"
$ws.Range("I5").Value2 = [double]"0.5257850852532363"
$ws.Range("J5").Value2 = [double]"0.5257850852532362"
$ws.Range("M5").Value2 = [double]"0.07271233333333334"
$ws.Range("N5").Value2 = [double]"0.218137"
$ws.Range("O5").Value2 = [double]"0.004171225362010892"
$ws.Range("P5").Value2 = [double]"0.004171225362010893"
$ws.Range("Q5").Value2 = [double]"0.6228802903852223"
$ws.Range("R5").Value2 = [double]"5.605922613467001"
$ws.Range("S5").Value2 = [double]"0.002193168082575359"
$ws.Range("T5").Value2 = [double]"0.002193168082575359"
$ws.Range("I6").Value2 = [double]"0.5257850852532363"
$ws.Range("J6").Value2 = [double]"0.5257850852532362"
$ws.Range("O6").Value2 = [double]"0.5387060579248023"
$ws.Range("P6").Value2 = [double]"0.5387060579248023"
$ws.Range("S6").Value2 = [double]"0.283243610592427"
$ws.Range("T6").Value2 = [double]"0.283243610592427"
$ws.Range("I7").Value2 = [double]"0.5257850852532363"
$ws.Range("J7").Value2 = [double]"0.5257850852532362"
$ws.Range("O7").Value2 = [double]"0.4571227167131868"
$ws.Range("P7").Value2 = [double]"0.4571227167131868"
$ws.Range("R7").Value2 = [double]"614.3505450677761"
$ws.Range("S7").Value2 = [double]"0.2403483065782339"
$ws.Range("T7").Value2 = [double]"0.2403483065782338"
$ws.Range("G8").Value2 = [double]"7.695814000000001"
$ws.Range("I8").Value2 = [double]"0.4723526081233437"
$ws.Range("J8").Value2 = [double]"0.4723526081233436"
$ws.Range("M8").Value2 = [double]"0.07271233333333334"
$ws.Range("N8").Value2 = [double]"0.218137"
$ws.Range("O8").Value2 = [double]"0.004171225362010892"
$ws.Range("P8").Value2 = [double]"0.004171225362010893"
$ws.Range("Q8").Value2 = [double]"0.5595805928393335"
$ws.Range("R8").Value2 = [double]"5.036225335554001"
$ws.Range("S8").Value2 = [double]"0.001970289178816084"
$ws.Range("T8").Value2 = [double]"0.001970289178816084"
$ws.Range("G9").Value2 = [double]"7.695814000000001"
$ws.Range("I9").Value2 = [double]"0.4723526081233437"
$ws.Range("J9").Value2 = [double]"0.4723526081233436"
$ws.Range("O9").Value2 = [double]"0.5387060579248023"
$ws.Range("P9").Value2 = [double]"0.5387060579248023"
$ws.Range("Q9").Value2 = [double]"72.26880091522467"
$ws.Range("R9").Value2 = [double]"650.4192082370221"
$ws.Range("S9").Value2 = [double]"0.2544592114726254"
$ws.Range("T9").Value2 = [double]"0.2544592114726254"
$ws.Range("G10").Value2 = [double]"7.695814000000001"
$ws.Range("I10").Value2 = [double]"0.4723526081233437"
$ws.Range("J10").Value2 = [double]"0.4723526081233436"
$ws.Range("O10").Value2 = [double]"0.4571227167131868"
$ws.Range("P10").Value2 = [double]"0.4571227167131868"
$ws.Range("Q10").Value2 = [double]"61.32418620876801"
$ws.Range("R10").Value2 = [double]"551.9176758789121"
$ws.Range("S10").Value2 = [double]"0.2159231074719022"
$ws.Range("T10").Value2 = [double]"0.2159231074719021"
